$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58. This pushes the former rows 58..106 down to 59..107,
# matching the target diff (the final row of the sheet becomes row 107).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with its data.
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 45126
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = "Tropicales y subtropicales"
$ws.Range("I58").Value = 100108003
$ws.Range("J58").Value = "Maracuyá"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 15
$ws.Range("N58").Value = 45000
$ws.Range("O58").Value = 45000
$ws.Range("P58").Value = 45000
$ws.Range("Q58").Value = '$/caja 18 kilos'
$ws.Range("R58").Value = "Región de Arica y Parinacota"
$ws.Range("S58").Value = 2500
$ws.Range("T58").Value = 18
